# Generate Report for Handback
# Updates the localization-status report:
#  - Marks the a18258af-de81-... row as failed ("Handback transform failed")
#    instead of "Ready for handoff" on every sheet that surfaces Status.
#  - Records the handback/handoff file-name mismatch error detail for both
#    the zh-cn and de-de locale sheets.
#  - Widens the "Error Detail" column (P) on the locale sheets so the new,
#    longer message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$failedStatus = "Handback transform failed"

# "Status" cells for the a18258af-de81-4dfc-82a3-839baf7ec22c.md row (row 3)
$wsOverview.Range("E3").Value = $failedStatus
$wsOverview.Range("F3").Value = $failedStatus
$wsZhCn.Range("C3").Value = $failedStatus
$wsDeDe.Range("C3").Value = $failedStatus

# New "Error Detail" messages describing the handback/handoff filename mismatch
$wsZhCn.Range("P3").Value = "Handback file name: pny311k3.1nq is different with handoff file name: a18258af-de81-4dfc-82a3-839baf7ec22c.d4fea3acfa0f7429ef8f42c47c8b1577c3899d8a.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: pny311k3.1nq is different with handoff file name: a18258af-de81-4dfc-82a3-839baf7ec22c.d4fea3acfa0f7429ef8f42c47c8b1577c3899d8a.de-de."

# Widen column P (Error Detail) on both locale sheets to fit the longer text
$wsZhCn.Columns(16).ColumnWidth = 39.17
$wsDeDe.Columns(16).ColumnWidth = 39.17
